$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$questionsText = @'
questions = [
    {
        "title": "A system administrator installed Elasticsearch on a server and ran the ./bin/elasticsearch -d -p pid -n 2  command to check if installation was successful. What will be the effect of this command?",
        "ques_type": 2,
        "options": [
            "It will start Elasticsearch on port 9100.",
            "It will start Elasticsearch in daemon mode showing the process ID of the service that is running.",
            "It will produce an error.",
            "It will start Elasticsearch in daemon mode with the number of nodes as two."
        ],
        "score": "It will produce an error."
    },
    {
        "title": "You want to create an Elasticsearch index using their REST APIs.What will happen if you run PUT /my-index-000001/_create?",
        "ques_type": 2,
        "options": [
            "It will create a new index named my-index-000001.",
            "It will return an error.",
            "It will create an index without any mapping.",
            "It will create a new index named index-000001."
        ],
        "score": "It will return an error."
    },
    {
        "title": "You have an Elasticsearch index named \u201canimals,\u201d which has fields including \u201chabitat,\u201d \u201cfood,\u201d and \u201cvertebrates.\u201d You want to get records where habitat or food have \u201csearch_string\u201d as a value.What code should you use?",
        "ques_type": 2,
        "options": [
            "GET animals/_search\n             {\n             \"query\": {\n             \"match\": {\n            \"query\" : \" search_string \"\n            , \"fields\": [\"habitat\",\"food\"]\n             }\n             }\n             }\n",
            "GET animals/_search\n             {\n             \"query\": {\n             \"match\": {\n            \"query\" : \" search \"\n            , \"fields\": [\"habitat\",\"food\"]\n             }\n             }\n             }\n",
            "GET animals/_search\n             {\n             \"query\": {\n             \"match\": {\n            \"query\" : \" search_string \"\n            , \"field\": \"food\"\n             },\n             \"match\": {\n            \"query\" : \"search\"\n            , \"field\": \"habitat\"\n             }\n             }\n             }\n",
            "GET animals/_search\n             {\n             \"query\": {\n             \"multi_match\": {\n            \"query\" : \" search_string \"\n            , \"fields\": [\"habitat\",\"food\"]\n             }\n             }\n             }\n"
        ],
        "score": "GET animals/_search\n             {\n             \"query\": {\n             \"multi_match\": {\n            \"query\" : \" search_string \"\n            , \"fields\": [\"habitat\",\"food\"]\n             }\n             }\n             }"
    },
    {
        "title": "You are using Kibana, an Elasticsearch query tool, as part of the application itself. What command should you use to hide this tool from the Kibana application?",
        "ques_type": 2,
        "options": [
            "elasticsearch.console.enabled: false",
            "console.ui.disabled: true",
            "console.enabled: false",
            "console.ui.enabled: false"
        ],
        "score": "console.ui.enabled: false"
    }
]
'@

$ws.Range("A2").ClearContents()
$ws.Range("A1").ClearFormats()
$ws.Range("A1").Value = $questionsText
$ws.Rows(1).AutoFit()
